# Update the SQL queries embedded in column B (TabQuery) and C (StatQuery)
# on the active sheet: the old schema used std.id / prt.id join keys;
# the new schema uses the fully-qualified study_id / participant_id keys.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsWithQueries = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cellsWithQueries) {
    $range = $ws.Range($addr)
    $text = $range.Text
    if ([string]::IsNullOrEmpty($text)) { continue }

    $new = $text.Replace(
        'df_participant prt ON std.id = prt."study.id"',
        'df_participant prt ON std.study_id = prt."study.study_id"')
    $new = $new.Replace(
        'df_diagnoses dgn ON prt.id = dgn."participant.id"',
        'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $new = $new.Replace(
        'df_treatments trt ON prt.id = trt."participant.id"',
        'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $new = $new.Replace(
        'df_treatment_resp trr ON prt.id = trr."participant.id"',
        'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $new = $new.Replace(
        'df_survival srv ON prt.id = srv."participant.id"',
        'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $new = $new.Replace(
        'df_reference_files rfs ON std.id = rfs."study.id"',
        'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

    if ($new -ne $text) {
        $range.Value = $new
    }
}

# Column C widened (and no longer marked as "best fit") to accommodate the
# longer join conditions in the updated queries.
$ws.Columns.Item(3).ColumnWidth = 71.43

Write-Host "Updated $($cellsWithQueries.Count) query cells and column C width"
